# "Generate Report for Archive"
#
# The CI report was regenerated: the localization status text changed from
# "Ready for handoff" to "In Translation" everywhere it's used (Overview!E:F
# and the per-locale Status column on zh-cn!C / de-de!C all point at the same
# shared string), and the now-shorter status text caused the report
# generator's column autofit to narrow the Status-ish columns.

$wb = $excel.ActiveWorkbook

# 1) Update the status text. All "Ready for handoff" cells share one string
#    table entry, so a single text reassignment (re-set on every cell that
#    shows it) flips every occurrence across all three sheets at once.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: compare with the literal on the LEFT. PowerShell's -eq
            # coerces the right operand to the left operand's type, and a
            # cell holding the text "True"/"False" reads back as a real
            # [bool] here -- "$trueBoolValue -eq 'anything non-empty'" would
            # itself coerce to $true and false-match every non-blank cell.
            if ("Ready for handoff" -eq $cell.Value()) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# 2) Narrow the columns that held the status text, matching the narrower
#    autofit width the report generator computed for the shorter string.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = 12.5   # F (de-de)

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # C (Status)
